$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-major fill of the new/changed data rows (2-9), columns A (Sending cluster) .. T (Edge total expression derived specificity)
# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Fzd4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.43424333333333
$ws.Range("H2").Value = 31.30273
$ws.Range("I2").Value = 0.9711091978791583
$ws.Range("J2").Value = 0.9711091978791584
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 19.48350866666667
$ws.Range("N2").Value = 58.450526
$ws.Range("O2").Value = 0.3081250754721727
$ws.Range("P2").Value = 0.3081250754721726
$ws.Range("Q2").Value = 203.2956704151089
$ws.Range("R2").Value = 1829.66103373598
$ws.Range("S2").Value = 0.2992230948882367
$ws.Range("T2").Value = 0.2992230948882367

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Fzd4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.43424333333333
$ws.Range("H3").Value = 31.30273
$ws.Range("I3").Value = 0.9711091978791583
$ws.Range("J3").Value = 0.9711091978791584
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 20.343383
$ws.Range("N3").Value = 61.03014900000001
$ws.Range("O3").Value = 0.3217236961512193
$ws.Range("P3").Value = 0.3217236961512193
$ws.Range("Q3").Value = 212.2678084451967
$ws.Range("R3").Value = 1910.41027600677
$ws.Range("S3").Value = 0.3124288405081286
$ws.Range("T3").Value = 0.3124288405081287

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Fzd4"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.43424333333333
$ws.Range("H4").Value = 31.30273
$ws.Range("I4").Value = 0.9711091978791583
$ws.Range("J4").Value = 0.9711091978791584
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1305583333333333
$ws.Range("N4").Value = 0.391675
$ws.Range("O4").Value = 0.002064735720865253
$ws.Range("P4").Value = 0.002064735720865253
$ws.Range("Q4").Value = 1.362277419194444
$ws.Range("R4").Value = 12.26049677275
$ws.Range("S4").Value = 0.002005083849721902
$ws.Range("T4").Value = 0.002005083849721902

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Fzd4"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 10.43424333333333
$ws.Range("H5").Value = 31.30273
$ws.Range("I5").Value = 0.9711091978791583
$ws.Range("J5").Value = 0.9711091978791584
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 23.275017
$ws.Range("N5").Value = 69.825051
$ws.Range("O5").Value = 0.3680864926557428
$ws.Range("P5").Value = 0.3680864926557428
$ws.Range("Q5").Value = 242.85719096547
$ws.Range("R5").Value = 2185.71471868923
$ws.Range("S5").Value = 0.3574521786330711
$ws.Range("T5").Value = 0.3574521786330711

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Fzd4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.310422
$ws.Range("H6").Value = 0.9312659999999999
$ws.Range("I6").Value = 0.02889080212084161
$ws.Range("J6").Value = 0.02889080212084161
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 19.48350866666667
$ws.Range("N6").Value = 58.450526
$ws.Range("O6").Value = 0.3081250754721727
$ws.Range("P6").Value = 0.3081250754721726
$ws.Range("Q6").Value = 6.048109727323999
$ws.Range("R6").Value = 54.43298754591599
$ws.Range("S6").Value = 0.008901980583935926
$ws.Range("T6").Value = 0.008901980583935926

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Fzd4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.310422
$ws.Range("H7").Value = 0.9312659999999999
$ws.Range("I7").Value = 0.02889080212084161
$ws.Range("J7").Value = 0.02889080212084161
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 20.343383
$ws.Range("N7").Value = 61.03014900000001
$ws.Range("O7").Value = 0.3217236961512193
$ws.Range("P7").Value = 0.3217236961512193
$ws.Range("Q7").Value = 6.315033637626001
$ws.Range("R7").Value = 56.83530273863401
$ws.Range("S7").Value = 0.009294855643090648
$ws.Range("T7").Value = 0.00929485564309065

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Wnt5a"
$ws.Range("C8").Value = "Fzd4"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.310422
$ws.Range("H8").Value = 0.9312659999999999
$ws.Range("I8").Value = 0.02889080212084161
$ws.Range("J8").Value = 0.02889080212084161
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1305583333333333
$ws.Range("N8").Value = 0.391675
$ws.Range("O8").Value = 0.002064735720865253
$ws.Range("P8").Value = 0.002064735720865253
$ws.Range("Q8").Value = 0.04052817894999999
$ws.Range("R8").Value = 0.3647536105499999
$ws.Range("S8").Value = 0.00005965187114335128
$ws.Range("T8").Value = 0.00005965187114335129

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Wnt5a"
$ws.Range("C9").Value = "Fzd4"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.310422
$ws.Range("H9").Value = 0.9312659999999999
$ws.Range("I9").Value = 0.02889080212084161
$ws.Range("J9").Value = 0.02889080212084161
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 23.275017
$ws.Range("N9").Value = 69.825051
$ws.Range("O9").Value = 0.3680864926557428
$ws.Range("P9").Value = 0.3680864926557428
$ws.Range("Q9").Value = 7.225077327174
$ws.Range("R9").Value = 65.025695944566
$ws.Range("S9").Value = 0.01063431402267168
$ws.Range("T9").Value = 0.01063431402267168

